$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 updates
$ws.Range("K2").Value = 3.85
$ws.Range("S2").Value = 2.98
$ws.Range("T2").Value = 1.66

# Row 4 updates
$ws.Range("F4").Value = 1.39
$ws.Range("G4").Value = 1.49

# Row 9 updates
$ws.Range("F9").Value = 1.82
$ws.Range("K9").Value = 5.6
